# "used more meaningful variables" -- rename the cryptic column headers
# Qk / Ak / ck on the "k" sheet to their descriptive equivalents, and
# update a few dependent data values / selections to match.

$wb = $excel.ActiveWorkbook

# --- Sheet "k": rename headers (Qk -> flowCapacity, Ak -> area, ck -> cost)
# Set D1 ("cost") before B1 ("flowCapacity") so the new shared-string table
# entries land in the same order as the target workbook (cost, then
# flowCapacity, appended after the existing unique strings).
$wsK = $wb.Worksheets.Item("k")
$wsK.Range("D1").Value = "cost"
$wsK.Range("B1").Value = "flowCapacity"
$wsK.Range("C1").Value = "area"
$wsK.Range("B2").Select() | Out-Null

# --- Sheet "misc": selection moved from B2 to B1
$wsMisc = $wb.Worksheets.Item("misc")
$wsMisc.Range("B1").Select() | Out-Null

# --- Sheet "m": updated flow-capacity figures for the first/second/third rows
# Activated/selected last so it ends up as the workbook's active tab.
$wsM = $wb.Worksheets.Item("m")
$wsM.Range("D2").Value = 300
$wsM.Range("D3").Value = 100
$wsM.Range("D4").Value = 300
$wsM.Range("E4").Value = 0.14
$wsM.Activate() | Out-Null
$wsM.Range("B3").Select() | Out-Null
